# Gantt Chart update: progress %, and re-sequenced task dates for the
# "Implementation and Testing" phase (rows 26-39), plus cursor/selection move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# --- Task progress updates -------------------------------------------------
# Row 9  "Complete Chart"  : 50% -> 100%
$ws.Range("D9").Value = 1
# Row 10 "Assign Tasks"    : 10% -> 20%
$ws.Range("D10").Value = 0.2

# --- Re-sequence the task dates for rows 26-39 ------------------------------
# Random room selection
$ws.Range("E26").Formula = "=E25"
$ws.Range("F26").Formula = "=E26+7"

# Random encounters
$ws.Range("E27").Formula = "=E25"
$ws.Range("F27").Formula = "=E27+14"

# Inventory
$ws.Range("E28").Formula = "=E25"
$ws.Range("F28").Formula = "=E28+14"

# Treasure/ Loot
$ws.Range("E29").Formula = "=F26"
$ws.Range("F29").Formula = "=E29+14"

# Game Over Screen ( Hp systerm)
$ws.Range("E30").Formula = "=F27"
$ws.Range("F30").Formula = "=E30+5"

# Scoreboard ( Enemy kills) - disucss
$ws.Range("E31").Formula = "=F30"
$ws.Range("F31").Formula = "=E31+3"

# Items
$ws.Range("E32").Formula = "=F27"
$ws.Range("F32").Formula = "=E32+4"

# Implementation and Testing (phase summary row)
$ws.Range("E33").Formula = "=MAX(F26:F32)"

# Test code
$ws.Range("E34").Formula = "=E33"
$ws.Range("F34").Formula = "=E34+7"

# Unit test
$ws.Range("E35").Formula = "=E33"
$ws.Range("F35").Formula = "=E35+5"

# Debugging
$ws.Range("E36").Formula = "=E34"
$ws.Range("F36").Formula = "=F34"

# Reflect on new ideas
$ws.Range("E37").Formula = "=E16"
$ws.Range("F37").Formula = "=E37+28"

# Optimise code
$ws.Range("E38").Formula = "=E36"
$ws.Range("F38").Formula = "=F36-14"

# Report
$ws.Range("E39").Formula = "=DATE(2026,1,6)-14"
$ws.Range("F39").Formula = "=DATE(2026,1,5)"

# --- Misc cosmetics: move the active selection/cursor -----------------------
$ws.Range("I8").Select()
